$wb = $excel.ActiveWorkbook

# --- Sheet: Neodymium ---
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C2").Value = [double]"0.0001624997844765925"
$ws.Range("D2").Value = [double]"0.4327962216169592"
$ws.Range("E2").Value = [double]"0.4954444338515565"

$ws.Range("B3").Value = [double]"2.183968937109617E-10"
$ws.Range("C3").Value = [double]"0.007863771111478995"
$ws.Range("D3").Value = [double]"0.3773862252567579"
$ws.Range("E3").Value = [double]"0.4405377924640967"

$ws.Range("B4").Value = [double]"3.409239115768776E-12"
$ws.Range("C4").Value = [double]"0.007112191862264317"
$ws.Range("D4").Value = [double]"0.3084289105169288"
$ws.Range("E4").Value = [double]"0.3890045279137931"

$ws.Range("C5").Value = [double]"1.58309364990403E-07"
$ws.Range("D5").Value = [double]"0.01705074686156489"
$ws.Range("E5").Value = [double]"0.03248552659786881"

# --- Sheet: Dysprosium ---
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C2").Value = [double]"0.000184130763029052"
$ws.Range("D2").Value = [double]"0.4904074105642924"
$ws.Range("E2").Value = [double]"0.5613949700759422"

$ws.Range("B3").Value = [double]"2.474685539534777E-10"
$ws.Range("C3").Value = [double]"0.008910548279841049"
$ws.Range("D3").Value = [double]"0.4276215740039337"
$ws.Range("E3").Value = [double]"0.4991794920271586"

$ws.Range("B4").Value = [double]"3.863056198855572E-12"
$ws.Range("C4").Value = [double]"0.00805892339258075"
$ws.Range("D4").Value = [double]"0.3494850828056441"
$ws.Range("E4").Value = [double]"0.4407864341311822"

$ws.Range("C5").Value = [double]"1.793825405013188E-07"
$ws.Range("D5").Value = [double]"0.01932043811594985"
$ws.Range("E5").Value = [double]"0.03680980143532343"

# --- Sheet: Copper ---
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("B2").Value = [double]"6.556944196948269E-06"
$ws.Range("C2").Value = [double]"0.005054593164513399"
$ws.Range("D2").Value = [double]"0.7294706924603505"
$ws.Range("E2").Value = [double]"0.6659852986845426"

$ws.Range("B3").Value = [double]"4.458494376266519E-05"
$ws.Range("C3").Value = [double]"0.01823510242611228"
$ws.Range("D3").Value = [double]"0.5180107211635228"
$ws.Range("E3").Value = [double]"0.5113972578016854"

$ws.Range("B4").Value = [double]"0.0001322251246896275"
$ws.Range("C4").Value = [double]"0.004881153221624548"
$ws.Range("D4").Value = [double]"0.4396160317803456"
$ws.Range("E4").Value = [double]"0.5150582967402831"

$ws.Range("B5").Value = [double]"4.153807974120016E-05"
$ws.Range("C5").Value = [double]"0.01070203110354613"
$ws.Range("D5").Value = [double]"0.6313641460866456"
$ws.Range("E5").Value = [double]"0.5201943099942047"

# --- Sheet: Raw silicon ---
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("B2").Value = [double]"4.966311329314387E-05"
$ws.Range("C2").Value = [double]"0.003454741537111926"
$ws.Range("D2").Value = [double]"0.9305532546032723"
$ws.Range("E2").Value = [double]"0.8621103429496831"

$ws.Range("B3").Value = [double]"5.299988190966854E-05"
$ws.Range("C3").Value = [double]"0.01154307604995803"
$ws.Range("D3").Value = [double]"0.4900770646459259"
$ws.Range("E3").Value = [double]"0.4790782953538905"

$ws.Range("B4").Value = [double]"0.0003396088080967769"
$ws.Range("C4").Value = [double]"0.003239926736874637"
$ws.Range("D4").Value = [double]"0.5064313401933364"
$ws.Range("E4").Value = [double]"0.6014344550262347"

$ws.Range("B5").Value = [double]"0.0001823578131693591"
$ws.Range("C5").Value = [double]"0.00411437408643171"
$ws.Range("D5").Value = [double]"0.8674230610438876"
$ws.Range("E5").Value = [double]"0.7152240119196929"
